$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B5/B6 ordering: B5 should be the "ARIA et accessibilité" text,
# B6 should be the new "Faire un site en utilisant SASS (projet 1)" text.
$ws.Range("B5").Value = "Faire un site respectant ARIA et accessibilité"
$ws.Range("B6").Value = "Faire un site en utilisant SASS (projet 1)"

# Add the 3 new rows (7, 8, 9) with dates and challenge text
$ws.Range("A7").Value = 44620
$ws.Range("B7").Value = "Faire un site en utilisant SASS + responsive (projet 1)"

$ws.Range("A8").Value = 44626
$ws.Range("B8").Value = "Utiliser des framework pour réaliser une landing page (Projet 2), un site e-commerce (Projet 3) et une boite mail (Projet 4)"

$ws.Range("A9").Value = 44626
$ws.Range("B9").Value = "Faire un projet pro pour un cv déployable sur LinkedIn "

# Apply same formatting as other date/challenge cells for the new rows
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("B8:B9").PasteSpecial(-4122)

# Update the selection to match the post-edit state
$ws.Range("J14").Select()
